$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing header values before overwriting (row 1, columns A:F)
$cantidad   = $ws.Range("A1").Value2
$titulo     = $ws.Range("B1").Value2
$descripcion= $ws.Range("C1").Value2
$precio     = $ws.Range("D1").Value2
$imagen     = $ws.Range("E1").Value2
$imagenLink = $ws.Range("F1").Value2

# Shift existing headers one column to the right to make room for the new
# "SKU" column, then write the new column in place (matching the new
# sharedStrings order: Cantidad, Titulo, Precio, Descripcion, Imagen_principal, link note, SKU)
$ws.Range("G1").Value = $imagenLink
$ws.Range("F1").Value = $imagen
$ws.Range("E1").Value = $precio
$ws.Range("D1").Value = $descripcion
$ws.Range("C1").Value = $titulo
$ws.Range("B1").Value = $cantidad
$ws.Range("A1").Value = "SKU"

$ws.Range("B4").Select()
